$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.735.33"
$ws.Range("E2").Value = "  -0.60%  "
$ws.Range("D3").Value = "3.862.94"
$ws.Range("E3").Value = "  +3.03%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'600.11"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.26%  "
$ws.Range("D6").Value = "'162.15"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.93%  "
$ws.Range("D7").Value = "3.859.25"
$ws.Range("E7").Value = "  +3.00%  "
$ws.Range("E9").Value = "  -1.75%  "
$ws.Range("D10").Value = "'0.167"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.76%  "
$ws.Range("E11").Value = "  -1.29%  "
$ws.Range("D12").Value = "'0.458"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.25%  "
$ws.Range("D13").Value = "'36.82"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.91%  "
$ws.Range("E14").Value = "  -2.05%  "
$ws.Range("D15").Value = "4.508.02"
$ws.Range("E15").Value = "  +3.03%  "
$ws.Range("D16").Value = "3.875.02"
$ws.Range("E16").Value = "  +2.92%  "
$ws.Range("D17").Value = "68.920.11"
$ws.Range("E17").Value = "  -0.33%  "
$ws.Range("D18").Value = "'7.56"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.86%  "
$ws.Range("E19").Value = "  -0.45%  "
$ws.Range("D20").Value = "'17.10"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.53%  "
$ws.Range("D21").Value = "'11.35"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.68%  "
$ws.Range("D22").Value = "'483.86"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.87%  "
$ws.Range("D24").Value = "'0.0000161"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +6.40%  "
$ws.Range("D25").Value = "'83.86"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.18%  "
$ws.Range("E26").Value = "  -2.87%  "
$ws.Range("D27").Value = "'12.07"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.69%  "
$ws.Range("E29").Value = "  -1.38%  "
$ws.Range("D30").Value = "'2.94"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.11%  "
$ws.Range("D31").Value = "4.013.83"
$ws.Range("E31").Value = "  +3.04%  "
$ws.Range("D32").Value = "'7.84"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.68%  "
$ws.Range("D33").Value = "'32.15"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.23%  "
$ws.Range("E34").Value = "  -4.25%  "
$ws.Range("D35").Value = "3.809.44"
$ws.Range("E35").Value = "  +3.44%  "
$ws.Range("E36").Value = "  -1.84%  "
$ws.Range("E37").Value = "  +1.78%  "
$ws.Range("E38").Value = "  +0.80%  "
$ws.Range("D39").Value = "'5.86"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.52%  "
$ws.Range("D40").Value = "'1.00"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.03%  "
$ws.Range("D41").Value = "'0.317"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.67%  "
$ws.Range("D42").Value = "'2.95"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.99%  "
$ws.Range("D43").Value = "'430.37"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.20%  "
$ws.Range("D44").Value = "'48.49"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.61%  "
$ws.Range("E45").Value = "  -0.58%  "
$ws.Range("E47").Value = "  -0.99%  "
$ws.Range("D48").Value = "'143.18"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.06%  "
$ws.Range("D49").Value = "2.838.79"
$ws.Range("E49").Value = "  +1.76%  "
$ws.Range("E50").Value = "  +1.06%  "
$ws.Range("D51").Value = "'25.78"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +12.81%  "
